# Fri, May 29, 2020  4:05:51 PM
# ---------------------------------------------------------------------------
# 1) Switch the deck's applied Design/Theme from "Integral" (Red Violet)
#    to the plain "Office Theme" palette. PowerPoint keeps the master's
#    theme part (theme1.xml) in place and just rewrites its twelve
#    scheme colors, which is exactly what editing the Theme Colors does.
# ---------------------------------------------------------------------------
$p = $ppt.ActivePresentation

$themeColors = $p.SlideMaster.Theme.ThemeColorScheme

$themeColors.Item(1).RGB  = "0x000000"   # Text/Background - Dark 1
$themeColors.Item(2).RGB  = "0xFFFFFF"   # Text/Background - Light 1
$themeColors.Item(3).RGB  = "0x44546A"   # Text/Background - Dark 2
$themeColors.Item(4).RGB  = "0xE7E6E6"   # Text/Background - Light 2
$themeColors.Item(5).RGB  = "0x5B9BD5"   # Accent 1
$themeColors.Item(6).RGB  = "0xED7D31"   # Accent 2
$themeColors.Item(7).RGB  = "0xA5A5A5"   # Accent 3
$themeColors.Item(8).RGB  = "0xFFC000"   # Accent 4
$themeColors.Item(9).RGB  = "0x4472C4"   # Accent 5
$themeColors.Item(10).RGB = "0x70AD47"   # Accent 6
$themeColors.Item(11).RGB = "0x0563C1"   # Hyperlink
$themeColors.Item(12).RGB = "0x954F72"   # Followed Hyperlink

# ---------------------------------------------------------------------------
# 2) The three tables that were left on the default/custom "Table_0" style
#    get switched over to the built-in table style that now matches the
#    new theme.
# ---------------------------------------------------------------------------
$newTableStyle = "{DD1F352F-4555-4026-8DDC-430E66866855}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    $tableShape = $slide.Shapes.Item(1)
    if ($tableShape.HasTable) {
        $tableShape.Table.ApplyStyle($newTableStyle)
    }
}
